$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 9 with the new "Week 8" lab report entry
# (order matches the shared-string insertion order: Date, Description, Time)
$ws.Range("A9").Value = "[Wk 8] Tuesday 1.5.18"
$ws.Range("D9").Value = "Laboratory report  2"
$ws.Range("B9").Value = "2200-0100"
$ws.Range("C9").Value = 3

# Recalculate so the SUBTOTAL in C21 reflects the new hours
$excel.Calculate()

# Move the active selection to D11, matching the author's final cursor position
$ws.Range("D11").Select()
